# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# These numbers were regenerated by the site's data pipeline (gh-pages output
# regenerated at commit 456a3b4), so only column F values change; everything
# else on the rows (including the "最低票价" column G) stays the same.

$wb = $excel.ActiveWorkbook

# Row -> (old, new) value for column F on sheet "展览"
$exhibitionUpdates = @{
    5  = 1721
    7  = 2175
    11 = 4886
    20 = 120
    21 = 3829
    22 = 701
    23 = 646
    27 = 115
    30 = 86
    34 = 925
    35 = 2437
    36 = 425
}

# Row -> new value for column F on sheet "全部类型"
$allTypesUpdates = @{
    5  = 1721
    7  = 2175
    11 = 4886
    20 = 120
    21 = 3829
    22 = 701
    23 = 646
    27 = 115
    30 = 86
    35 = 925
    36 = 2437
    37 = 425
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
